# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: D.Prescott
$rushing.Range("C2").Value = 13
$rushing.Range("D2").Value = 6
$rushing.Range("E2").Value = 14
$rushing.Range("F2").Value = 11

# Row 4: E.Elliott
$rushing.Range("C4").Value = 131
$rushing.Range("D4").Value = 61
$rushing.Range("E4").Value = 28

# Row 5: T.Pollard
$rushing.Range("C5").Value = 77
$rushing.Range("D5").Value = 44

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: E.Elliott
$receiving.Range("C2").Value = 47
$receiving.Range("D2").Value = 34

# Row 3: T.Pollard
$receiving.Range("C3").Value = 36
$receiving.Range("D3").Value = 29

# Row 5: A.Cooper
$receiving.Range("C5").Value = 73
$receiving.Range("D5").Value = 51
$receiving.Range("G5").Value = 19
$receiving.Range("H5").Value = 15

# Row 6: C.Wilson
$receiving.Range("C6").Value = 79
$receiving.Range("D6").Value = 55
$receiving.Range("E6").Value = 35
$receiving.Range("F6").Value = 19

# Row 7: M.Gallup
$receiving.Range("C7").Value = 47
$receiving.Range("D7").Value = 36
$receiving.Range("E7").Value = 15
$receiving.Range("F7").Value = 9

# Row 8: N.Brown
$receiving.Range("C8").Value = 29
$receiving.Range("D8").Value = 24
$receiving.Range("E8").Value = 11
$receiving.Range("F8").Value = 7
$receiving.Range("G8").Value = 6
$receiving.Range("H8").Value = 3

# Row 12: D.Schultz
$receiving.Range("C12").Value = 76
$receiving.Range("D12").Value = 61
$receiving.Range("E12").Value = 10
$receiving.Range("G12").Value = 12
